$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Cluster Names"
$ws.Range("B1").Value = "Active Cases"

$names = @(
    "3323 Villa Maria Catholic Homes St Bernadette's Aged Care Sunshine North",
    "3398 BlueCross Elly Kay Mordialloc",
    "3601 Baptcare Westhaven community",
    "3653 Fronditha Thalpori St Albans Aged Care",
    "3975 Aurrum Aged Care Brunswick West",
    "3988 Kerala Manor Aged Care Diamond Creek",
    "4257 BlueCross The Gables Camberwell",
    "4295 Hope Aged Care Sunshine West",
    "44087 Fitzroy Primary School Fitzroy",
    "44098 Stawell Primary School",
    "44234 Lucknow Primary School Bairnsdale",
    "44366 Lysterfield Primary School Lysterfield",
    "44444 Nar Nar Goon Primary School Nar Nar Goon",
    "44630 Black Rock Primary School Black Rock",
    "44666 Gardenvale Primary School Senior School Campus Brighton East",
    "44811 Dandenong North Primary School Dandenong",
    "44865 Parktone Primary School Parkdale",
    "44950 Templestowe Valley Primary School Templestowe Lower",
    "44982 Diamond Creek East Primary School Diamond Creek",
    "45026 Churchill North Primary School Churchill",
    "45248 Brookside P-9 College Caroline Springs",
    "45249 Creekside K-9 College Caroline Springs",
    "45267 Epping Views Primary School Epping",
    "45315 Red Hill Consolidated School Red Hill",
    "45518 Ashwood High School Ashwood",
    "45569 Nhill College Nhill",
    "45585 Mount Ridley College Craigieburn",
    "45648 St Brendans Primary School Shepparton",
    "4574 Village Glen Aged Care Residences Mornington",
    "45755 St Patricks Catholic Parish Primary School Mentone",
    "45784 Holy Rosary Primary School White Hills",
    "45846 St Mary's School Mooroopna Outbreak",
    "45848 St Kevin's College Toorak Glendalough Campus Junior School",
    "45912 St Bernadette's Catholic Primary School Sunshine North",
    "45950 St Luke's Primary School Lalor",
    "46028 St Anne's Catholic Primary School Sunbury",
    "46037 Nazareth Catholic Primary School Grovedale",
    "46050 Our Lady's Catholic Primary School Craigieburn",
    "46052 St. Francis of Assisi Primary School Mill Park",
    "46093 St Brendan's Primary School Somerville",
    "46095 Bethany Catholic Primary School Werribee",
    "46105 Christ the Priest Primary School Caroline Springs",
    "46117 Marymede Catholic College South Morang",
    "46125 Our Lady of the Southern Cross Primary School Manor Lakes",
    "46221 Bialik College Hawthorn Outbreak Site",
    "46239 Gilson College Taylors Hill",
    "46287 Oakleigh Grammar Melbourne Private School Oakleigh",
    "46390 Al Siraat College Epping",
    "50584 St Mary of the Cross MacKillop Primary School Epping",
    "51529 Sirius College Primary School Dallas",
    "Alfred Health The Alfred Hospital Melbourne",
    "Covenant College Bell Post Hill",
    "Epping Views Primary School Camp Cape Schanck",
    "Hamilton Country Music Festival Hamilton Golf Club Hamilton",
    "House Party 27 November Private Residence Brunswick West",
    "Islamic College of Melbourne Tarneit",
    "Little Munchkins Childcare Centre Hillside",
    "Social Gathering 20 November Sunbury",
    "Springside Primary School Caroline Springs",
    "St Josephs Catholic Primary School Warragul",
    "St Vincents Hospital Melbourne Emergency Department Fitzroy",
    "The Village Early Learning Centre Sandringham",
    "Wagstaff Meat Processing Plant Cranbourne East",
    "Western Health Sunshine Hospital Emergency Department St Albans"
)
$values = @(14, 37, 20, 22, 10, 10, 28, 32, 22, 24, 13, 11, 18, 21, 11, 20, 22, 65, 11, 11, 30, 12, 20, 13, 21, 33, 10, 28, 11, 13, 26, 15, 16, 10, 17, 12, 27, 12, 26, 14, 11, 40, 13, 37, 13, 12, 25, 30, 13, 14, 13, 23, 14, 12, 10, 16, 10, 20, 23, 15, 12, 13, 36, 11)

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $names[$i]
    $ws.Cells.Item($r, 2).Value = $values[$i]
}

Write-Host ("rows written: " + $names.Length)